$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 is treated as text so numeric-looking strings (e.g. "561.50")
# are not auto-converted to numbers by Excel, matching the original inlineStr cells.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.532.02"
$ws.Range("D3").Value = "3.361.21"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "561.50"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "176.40"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").Value = "3.351.47"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.632"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  +8.00%  "
$ws.Range("D12").Value = "55.10"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "9.10"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "3.905.12"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "18.30"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "3.364.77"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "11.84"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "64.512.80"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "463.76"
$ws.Range("E22").Value = "  +13.40%  "
$ws.Range("D23").Value = "4.84"
$ws.Range("E23").Value = "  +9.77%  "
$ws.Range("D24").Value = "4.11"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "86.34"
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("D26").Value = "13.49"
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("D29").Value = "8.80"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Value = "30.17"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "6.72"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "11.50"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "579.17"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "59.41"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "0.140"
$ws.Range("E37").Value = "  -6.82%  "
$ws.Range("D38").Value = "35.90"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.47"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  +4.29%  "
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "3.094.08"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "8.36"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "136.00"
$ws.Range("E51").Value = "  +1.44%  "

# Restore the default (unstyled) cell style now that values are set as text,
# so the saved file does not carry a spurious explicit style index.
$numRange.Style = "Normal"

